$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 9).Value = 'sd'
$ws.Cells.Item(6, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(8, 9).Value = 'aa'
$ws.Cells.Item(8, 10).Value = 'Agree/Accept'
$ws.Cells.Item(9, 9).Value = 'sd'
$ws.Cells.Item(9, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(15, 9).Value = 'sd'
$ws.Cells.Item(15, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(34, 9).Value = 'sv'
$ws.Cells.Item(34, 10).Value = 'Statement-opinion'
$ws.Cells.Item(41, 9).Value = 'sd'
$ws.Cells.Item(41, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(58, 9).Value = 'sv'
$ws.Cells.Item(58, 10).Value = 'Statement-opinion'
$ws.Cells.Item(61, 9).Value = 'aa'
$ws.Cells.Item(61, 10).Value = 'Agree/Accept'
$ws.Cells.Item(88, 9).Value = 'sv'
$ws.Cells.Item(88, 10).Value = 'Statement-opinion'
$ws.Cells.Item(99, 9).Value = 'sd'
$ws.Cells.Item(99, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(103, 9).Value = 'aa'
$ws.Cells.Item(103, 10).Value = 'Agree/Accept'
$ws.Cells.Item(157, 9).Value = 'sd'
$ws.Cells.Item(157, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(159, 9).Value = 'sv'
$ws.Cells.Item(159, 10).Value = 'Statement-opinion'
$ws.Cells.Item(172, 9).Value = 'sd'
$ws.Cells.Item(172, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(175, 9).Value = 'sv'
$ws.Cells.Item(175, 10).Value = 'Statement-opinion'
$ws.Cells.Item(181, 9).Value = 'sv'
$ws.Cells.Item(181, 10).Value = 'Statement-opinion'
$ws.Cells.Item(186, 9).Value = 'sd'
$ws.Cells.Item(186, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(188, 9).Value = 'sd'
$ws.Cells.Item(188, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(189, 9).Value = 'sv'
$ws.Cells.Item(189, 10).Value = 'Statement-opinion'
$ws.Cells.Item(193, 9).Value = 'sd'
$ws.Cells.Item(193, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(194, 9).Value = 'sv'
$ws.Cells.Item(194, 10).Value = 'Statement-opinion'
$ws.Cells.Item(197, 9).Value = 'sd'
$ws.Cells.Item(197, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(199, 9).Value = 'sd'
$ws.Cells.Item(199, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(201, 9).Value = 'sv'
$ws.Cells.Item(201, 10).Value = 'Statement-opinion'
$ws.Cells.Item(204, 9).Value = 'sd'
$ws.Cells.Item(204, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(213, 9).Value = 'sv'
$ws.Cells.Item(213, 10).Value = 'Statement-opinion'
$ws.Cells.Item(219, 9).Value = 'sd'
$ws.Cells.Item(219, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(222, 9).Value = 'sv'
$ws.Cells.Item(222, 10).Value = 'Statement-opinion'
$ws.Cells.Item(230, 9).Value = 'sv'
$ws.Cells.Item(230, 10).Value = 'Statement-opinion'
$ws.Cells.Item(231, 9).Value = 'b'
$ws.Cells.Item(231, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(241, 9).Value = 'aa'
$ws.Cells.Item(241, 10).Value = 'Agree/Accept'
$ws.Cells.Item(243, 9).Value = 'sd'
$ws.Cells.Item(243, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(254, 9).Value = 'ba'
$ws.Cells.Item(254, 10).Value = 'Appreciation'
$ws.Cells.Item(259, 9).Value = 'sv'
$ws.Cells.Item(259, 10).Value = 'Statement-opinion'
$ws.Cells.Item(265, 9).Value = 'sd'
$ws.Cells.Item(265, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(269, 9).Value = 'qy'
$ws.Cells.Item(269, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(281, 9).Value = 'sv'
$ws.Cells.Item(281, 10).Value = 'Statement-opinion'
$ws.Cells.Item(285, 9).Value = 'sv'
$ws.Cells.Item(285, 10).Value = 'Statement-opinion'
$ws.Cells.Item(305, 9).Value = 'qy'
$ws.Cells.Item(305, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(321, 9).Value = 'sd'
$ws.Cells.Item(321, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(323, 9).Value = 'aa'
$ws.Cells.Item(323, 10).Value = 'Agree/Accept'
$ws.Cells.Item(324, 9).Value = 'aa'
$ws.Cells.Item(324, 10).Value = 'Agree/Accept'
$ws.Cells.Item(343, 9).Value = 'ba'
$ws.Cells.Item(343, 10).Value = 'Appreciation'
$ws.Cells.Item(348, 9).Value = 'ba'
$ws.Cells.Item(348, 10).Value = 'Appreciation'
$ws.Cells.Item(351, 9).Value = 'sd'
$ws.Cells.Item(351, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(357, 9).Value = 'sd'
$ws.Cells.Item(357, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(364, 9).Value = '%'
$ws.Cells.Item(364, 10).Value = 'Uninterpretable'
$ws.Cells.Item(367, 9).Value = 'aa'
$ws.Cells.Item(367, 10).Value = 'Agree/Accept'
$ws.Cells.Item(369, 9).Value = 'sd'
$ws.Cells.Item(369, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(383, 9).Value = '%'
$ws.Cells.Item(383, 10).Value = 'Uninterpretable'
$ws.Cells.Item(385, 9).Value = 'aa'
$ws.Cells.Item(385, 10).Value = 'Agree/Accept'
$ws.Cells.Item(388, 9).Value = 'aa'
$ws.Cells.Item(388, 10).Value = 'Agree/Accept'
$ws.Cells.Item(390, 9).Value = 'sd'
$ws.Cells.Item(390, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(401, 9).Value = 'sv'
$ws.Cells.Item(401, 10).Value = 'Statement-opinion'
$ws.Cells.Item(409, 9).Value = 'ba'
$ws.Cells.Item(409, 10).Value = 'Appreciation'
$ws.Cells.Item(416, 9).Value = 'sd'
$ws.Cells.Item(416, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(420, 9).Value = 'sv'
$ws.Cells.Item(420, 10).Value = 'Statement-opinion'
$ws.Cells.Item(422, 9).Value = 'sv'
$ws.Cells.Item(422, 10).Value = 'Statement-opinion'
$ws.Cells.Item(440, 9).Value = 'aa'
$ws.Cells.Item(440, 10).Value = 'Agree/Accept'
$ws.Cells.Item(457, 9).Value = 'b'
$ws.Cells.Item(457, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(459, 9).Value = 'ba'
$ws.Cells.Item(459, 10).Value = 'Appreciation'
$ws.Cells.Item(460, 9).Value = 'sd'
$ws.Cells.Item(460, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(475, 9).Value = 'ba'
$ws.Cells.Item(475, 10).Value = 'Appreciation'
$ws.Cells.Item(494, 9).Value = 'aa'
$ws.Cells.Item(494, 10).Value = 'Agree/Accept'
$ws.Cells.Item(498, 9).Value = 'aa'
$ws.Cells.Item(498, 10).Value = 'Agree/Accept'
$ws.Cells.Item(499, 9).Value = 'sd'
$ws.Cells.Item(499, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(501, 9).Value = 'ba'
$ws.Cells.Item(501, 10).Value = 'Appreciation'
$ws.Cells.Item(508, 9).Value = 'sd'
$ws.Cells.Item(508, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(513, 9).Value = 'sv'
$ws.Cells.Item(513, 10).Value = 'Statement-opinion'
$ws.Cells.Item(516, 9).Value = '%'
$ws.Cells.Item(516, 10).Value = 'Uninterpretable'
$ws.Cells.Item(521, 9).Value = 'sv'
$ws.Cells.Item(521, 10).Value = 'Statement-opinion'
$ws.Cells.Item(524, 9).Value = 'qy'
$ws.Cells.Item(524, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(525, 9).Value = 'qy'
$ws.Cells.Item(525, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(530, 9).Value = 'sv'
$ws.Cells.Item(530, 10).Value = 'Statement-opinion'
$ws.Cells.Item(542, 9).Value = 'ba'
$ws.Cells.Item(542, 10).Value = 'Appreciation'
$ws.Cells.Item(544, 9).Value = 'sd'
$ws.Cells.Item(544, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(546, 9).Value = 'aa'
$ws.Cells.Item(546, 10).Value = 'Agree/Accept'
$ws.Cells.Item(553, 9).Value = 'b'
$ws.Cells.Item(553, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(559, 9).Value = '%'
$ws.Cells.Item(559, 10).Value = 'Uninterpretable'
$ws.Cells.Item(560, 9).Value = 'b'
$ws.Cells.Item(560, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(572, 9).Value = 'sv'
$ws.Cells.Item(572, 10).Value = 'Statement-opinion'
$ws.Cells.Item(573, 9).Value = 'sv'
$ws.Cells.Item(573, 10).Value = 'Statement-opinion'
$ws.Cells.Item(578, 9).Value = 'sv'
$ws.Cells.Item(578, 10).Value = 'Statement-opinion'
$ws.Cells.Item(579, 9).Value = 'aa'
$ws.Cells.Item(579, 10).Value = 'Agree/Accept'
$ws.Cells.Item(581, 9).Value = 'ba'
$ws.Cells.Item(581, 10).Value = 'Appreciation'
$ws.Cells.Item(596, 9).Value = '%'
$ws.Cells.Item(596, 10).Value = 'Uninterpretable'
$ws.Cells.Item(606, 9).Value = 'b'
$ws.Cells.Item(606, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(615, 9).Value = 'ba'
$ws.Cells.Item(615, 10).Value = 'Appreciation'
